# Baltimore GDP data refresh ("updating GDP and algo")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Revised GDP figures for existing years (rows 12-30) ---
$ws.Cells.Item(12,2).Value = 102974.67200000001
$ws.Cells.Item(13,2).Value = 108228.454
$ws.Cells.Item(14,2).Value = 112798.30499999999
$ws.Cells.Item(15,2).Value = 119753.76700000001
$ws.Cells.Item(16,2).Value = 128559.908
$ws.Cells.Item(17,2).Value = 135614.829
$ws.Cells.Item(18,2).Value = 140017.144
$ws.Cells.Item(19,2).Value = 144649.59400000001
$ws.Cells.Item(20,2).Value = 146563.97
$ws.Cells.Item(21,2).Value = 153461.467
$ws.Cells.Item(22,2).Value = 158810.78899999999
$ws.Cells.Item(23,2).Value = 163424.44500000001
$ws.Cells.Item(24,2).Value = 169217.63500000001
$ws.Cells.Item(25,2).Value = 175334.84400000001
$ws.Cells.Item(26,2).Value = 183664.932
$ws.Cells.Item(27,2).Value = 192223.283
$ws.Cells.Item(28,2).Value = 198751.55900000001
$ws.Cells.Item(29,2).Value = 205653.04699999999
$ws.Cells.Item(30,2).Value = 212887.514

# --- New observation: 2020-01-01 (row 31) ---
$ws.Cells.Item(31,1).Value = 43831
$ws.Cells.Item(31,1).NumberFormat = "yyyy\-mm\-dd"
$ws.Cells.Item(31,2).Value = 205810.92300000001
$ws.Cells.Item(31,2).NumberFormat = "0.000"

# --- Selection state left by the editing user ---
$ws.Range("A1:B1048576").Select()
